$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "46.597.36"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.611.61"
$ws.Range("E3").Value = "  +7.77%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.13"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.89"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("E7").Value = "  +5.98%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  +14.04%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.55"
$ws.Range("E10").Value = "  +11.80%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +8.05%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.49"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.20"
$ws.Range("E13").Value = "  +14.88%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.003.53"
$ws.Range("E14").Value = "  +7.35%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.106"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.606.40"
$ws.Range("E16").Value = "  +6.98%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.928"
$ws.Range("E17").Value = "  +9.35%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.10"
$ws.Range("E18").Value = "  +7.46%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "46.729.34"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("E20").Value = "  +7.62%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.82"
$ws.Range("E22").Value = "  +9.26%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.93"
$ws.Range("E23").Value = "  +13.53%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.89"
$ws.Range("E24").Value = "  +6.56%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.07"
$ws.Range("E25").Value = "  +9.38%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  +11.81%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.16"
$ws.Range("E27").Value = "  +36.34%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.04"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.72"
$ws.Range("E30").Value = "  +9.58%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.32"
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "39.41"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.39"
$ws.Range("E33").Value = "  +15.21%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.64"
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0845"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.84"
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").Value = "  +11.82%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.22"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("E39").Value = "  +9.65%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.123"
$ws.Range("E40").Value = "  +6.24%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.35"
$ws.Range("E41").Value = "  +42.45%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.29"
$ws.Range("E42").Value = "  +6.49%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0335"
$ws.Range("E43").Value = "  +10.95%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  +12.74%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.10"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.140.24"
$ws.Range("E46").Value = "  +8.21%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.35"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.58"
$ws.Range("E49").Value = "  +11.17%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.80"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.05"
$ws.Range("E51").Value = "  +9.04%  "
